# Stomp Master Data.xlsx - apply "Parse excel file and produce sql script" edit
$wb = $excel.ActiveWorkbook

# --- Teachers sheet: add a new teacher row (George Write / Columbus) ---
$wsTeachers = $wb.Worksheets.Item("Teachers")
$wsTeachers.Range("A3").Value = "George"
$wsTeachers.Range("B3").Value = "Write"
$wsTeachers.Range("C3").Value = "george.wright@columbus.edu"
$wsTeachers.Range("D3").Value = "Columbus"
$wsTeachers.Range("C3").Hyperlinks.Add($wsTeachers.Range("C3"), "mailto:george.wright@columbus.edu") | Out-Null
$wsTeachers.Range("C2").Copy() | Out-Null
$wsTeachers.Range("C3").PasteSpecial(-4122) | Out-Null
$wsTeachers.Range("C2").Select() | Out-Null

# --- Stompers sheet: add headers + two fellows (Sam Heilbron, Terrence Roh) ---
$wsStompers = $wb.Worksheets.Item("Stompers")
$wsStompers.Range("A1").Value = "First Name"
$wsStompers.Range("B1").Value = "Last Name"
$wsStompers.Range("C1").Value = "Email"
$wsStompers.Range("D1").Value = "UTLN"

$wsStompers.Range("A2").Value = "Sam"
$wsStompers.Range("B2").Value = "Heilbron"
$wsStompers.Range("C2").Value = "samheilbron@gmail.com"
$wsStompers.Range("D2").Value = "sheilb01"

$wsStompers.Range("A3").Value = "Terrence"
$wsStompers.Range("B3").Value = "Roh"
$wsStompers.Range("C3").Value = "troh@gmail.com"
$wsStompers.Range("D3").Value = "troh01"

$wsStompers.Range("C2").Hyperlinks.Add($wsStompers.Range("C2"), "mailto:samheilbron@gmail.com") | Out-Null
$wsStompers.Range("C3").Hyperlinks.Add($wsStompers.Range("C3"), "mailto:troh@gmail.com") | Out-Null

$wsStompers.Range("A1:D1").Font.Bold = $true
$wsStompers.Columns.Item(1).ColumnWidth = 10.83203125
$wsStompers.Columns.Item(2).ColumnWidth = 10.83203125
$wsStompers.Columns.Item(3).ColumnWidth = 21.6640625
$wsStompers.Columns.Item(4).ColumnWidth = 10.83203125
$wsStompers.Range("C15").Select() | Out-Null

# --- Pairings sheet: update the first two pairing rows ---
$wsPairings = $wb.Worksheets.Item("Pairings")
$wsPairings.Range("E2").Value = "anna.mcCormick@brooks.edu"
$wsPairings.Range("F2").Value = "Sam Heilbron"
$wsPairings.Range("G2").Value = "Terrence Roh"

$wsPairings.Range("D3").Value = "10:20-11:45"
$wsPairings.Range("E3").Value = "george.wright@columbus.edu"
$wsPairings.Range("F3").Value = "Sam Heilbron"
$wsPairings.Range("G3").Value = "Terrence Roh"

$wsPairings.Range("E3").Hyperlinks.Add($wsPairings.Range("E3"), "mailto:george.wright@columbus.edu") | Out-Null
$wsPairings.Range("E2").Hyperlinks.Add($wsPairings.Range("E2"), "mailto:anna.mcCormick@brooks.edu") | Out-Null
$wsPairings.Columns.Item(5).ColumnWidth = 32.83203125
$wsPairings.Range("E2").Select() | Out-Null

# --- Activate Pairings as the visible tab (moves tabSelected / activeTab) ---
$wsPairings.Activate()
